$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range($cellRef).Formula = '=""&"' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue 'D2' '45.350.51'
Set-TextValue 'E2' '  +5.21%  '
Set-TextValue 'D3' '2.456.98'
Set-TextValue 'E3' '  +3.63%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.07%  '
Set-TextValue 'D5' '319.83'
Set-TextValue 'E5' '  +5.35%  '
Set-TextValue 'D6' '104.50'
Set-TextValue 'E6' '  +8.54%  '
Set-TextValue 'E7' '  +2.80%  '
Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  -0.04%  '
Set-TextValue 'D9' '0.532'
Set-TextValue 'E9' '  +10.24%  '
Set-TextValue 'D10' '36.02'
Set-TextValue 'E10' '  +4.68%  '
Set-TextValue 'D11' '0.0806'
Set-TextValue 'E11' '  +2.09%  '
Set-TextValue 'E12' '  -2.48%  '
Set-TextValue 'D13' '18.56'
Set-TextValue 'E13' '  +1.31%  '
Set-TextValue 'E14' '  +3.83%  '
Set-TextValue 'D15' '2.841.16'
Set-TextValue 'E15' '  +3.81%  '
Set-TextValue 'D16' '2.460.14'
Set-TextValue 'E16' '  +1.37%  '
Set-TextValue 'E17' '  +4.83%  '
Set-TextValue 'D18' '45.255.69'
Set-TextValue 'E18' '  +5.04%  '
Set-TextValue 'D19' '12.39'
Set-TextValue 'E19' '  +3.27%  '
Set-TextValue 'E20' '  +1.36%  '
Set-TextValue 'E21' '  +4.57%  '
Set-TextValue 'D22' '69.26'
Set-TextValue 'E22' '  +1.79%  '
Set-TextValue 'D23' '244.46'
Set-TextValue 'E23' '  +3.75%  '
Set-TextValue 'E24' '  +3.36%  '
Set-TextValue 'D25' '2.52'
Set-TextValue 'E25' '  +3.15%  '
Set-TextValue 'E26' '  +0.07%  '
Set-TextValue 'D27' '25.55'
Set-TextValue 'E27' '  +4.31%  '
Set-TextValue 'D28' '2.19'
Set-TextValue 'E28' '  -7.58%  '
Set-TextValue 'D29' '9.60'
Set-TextValue 'E29' '  +2.45%  '
Set-TextValue 'D30' '34.03'
Set-TextValue 'E30' '  +6.76%  '
Set-TextValue 'D31' '49.64'
Set-TextValue 'E31' '  +3.34%  '
Set-TextValue 'E32' '  +15.29%  '
Set-TextValue 'D33' '20.44'
Set-TextValue 'E33' '  +14.44%  '
Set-TextValue 'D34' '5.26'
Set-TextValue 'E34' '  +4.12%  '
Set-TextValue 'E35' '  +0.16%  '
Set-TextValue 'D36' '0.0767'
Set-TextValue 'E36' '  +3.99%  '
Set-TextValue 'D37' '1.92'
Set-TextValue 'E37' '  +5.36%  '
Set-TextValue 'D38' '4.54'
Set-TextValue 'E38' '  +4.58%  '
Set-TextValue 'E39' '  +1.03%  '
Set-TextValue 'D40' '125.30'
Set-TextValue 'E40' '  -2.20%  '
Set-TextValue 'E41' '  +2.40%  '
Set-TextValue 'E42' '  -2.85%  '
Set-TextValue 'D43' '21.35'
Set-TextValue 'E43' '  +0.61%  '
Set-TextValue 'E44' '  +4.76%  '
Set-TextValue 'D45' '1.950.81'
Set-TextValue 'E45' '  +1.12%  '
Set-TextValue 'D46' '3.00'
Set-TextValue 'E46' '  +8.48%  '
Set-TextValue 'E47' '  -0.81%  '
Set-TextValue 'D48' '9.27'
Set-TextValue 'E48' '  +0.31%  '
Set-TextValue 'D49' '1.79'
Set-TextValue 'E49' '  +17.19%  '
Set-TextValue 'D50' '76.43'
Set-TextValue 'E50' '  +6.58%  '
Set-TextValue 'D51' '54.02'
Set-TextValue 'E51' '  +4.31%  '

$excel.CutCopyMode = 0

